# Fix to slide 8: move the two topology diagrams (component boxes, ovals,
# connectors) to the right, and reposition/resize the explanatory text box.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

$sh = $s.Shapes.Item("Rounded Rectangle 3")
$sh.Left = 271.39503937007873
$sh.Top = 96.10811023622047

$sh = $s.Shapes.Item("TextBox 4")
$sh.Left = 291.9134645669291
$sh.Top = 100.32582677165354

$sh = $s.Shapes.Item("Oval 5")
$sh.Left = 292.5425196850394
$sh.Top = 135.91299212598426

$sh = $s.Shapes.Item("Oval 6")
$sh.Left = 350.7930708661417
$sh.Top = 124.05062992125984

$sh = $s.Shapes.Item("Oval 7")
$sh.Left = 366.8511023622047
$sh.Top = 166.7551968503937

$sh = $s.Shapes.Item("Straight Connector 8")
$sh.Left = 312.61503937007876
$sh.Top = 135.91299212598426

$sh = $s.Shapes.Item("Straight Arrow Connector 9")
$sh.Left = 360.82929133858266
$sh.Top = 147.77543307086614

$sh = $s.Shapes.Item("Straight Arrow Connector 10")
$sh.Left = 309.6755118110236
$sh.Top = 156.16338582677164

$sh = $s.Shapes.Item("Curved Connector 11")
$sh.Left = 354.2525196850394
$sh.Top = 130.62748031496062

$sh = $s.Shapes.Item("Rectangle 12")
$sh.Left = 292.5425196850394
$sh.Top = 205.24212598425197

$sh = $s.Shapes.Item("Rectangle 13")
$sh.Left = 390.9381102362205
$sh.Top = 205.76937007874017

$sh = $s.Shapes.Item("TextBox 14")
$sh.Left = 310.72779527559055
$sh.Top = 207.58244094488188

$sh = $s.Shapes.Item("TextBox 15")
$sh.Left = 407.88826771653544
$sh.Top = 207.0552755905512

$sh = $s.Shapes.Item("Rounded Rectangle 16")
$sh.Left = 525.2010236220473
$sh.Top = 104.93905511811023

$sh = $s.Shapes.Item("TextBox 17")
$sh.Left = 546.7255118110236
$sh.Top = 109.1567716535433

$sh = $s.Shapes.Item("Oval 18")
$sh.Left = 546.3485039370079
$sh.Top = 144.74393700787402

$sh = $s.Shapes.Item("Oval 19")
$sh.Left = 604.5989763779528
$sh.Top = 132.8814960629921

$sh = $s.Shapes.Item("Oval 20")
$sh.Left = 620.6570078740158
$sh.Top = 175.58614173228347

$sh = $s.Shapes.Item("Straight Connector 21")
$sh.Left = 566.4210236220473
$sh.Top = 144.74393700787402

$sh = $s.Shapes.Item("Straight Arrow Connector 22")
$sh.Left = 614.6352755905511
$sh.Top = 156.60629921259843

$sh = $s.Shapes.Item("Straight Arrow Connector 23")
$sh.Left = 563.4814960629922
$sh.Top = 164.9943307086614

$sh = $s.Shapes.Item("Curved Connector 24")
$sh.Left = 608.0584251968504
$sh.Top = 139.4583464566929

$sh = $s.Shapes.Item("Rectangle 25")
$sh.Left = 546.3485039370079
$sh.Top = 214.07307086614173

$sh = $s.Shapes.Item("Rectangle 26")
$sh.Left = 644.7440944881889
$sh.Top = 214.60031496062993

$sh = $s.Shapes.Item("TextBox 27")
$sh.Left = 562.3261417322834
$sh.Top = 218.73007874015747

$sh = $s.Shapes.Item("TextBox 28")
$sh.Left = 660.8020472440945
$sh.Top = 218.02708661417321

$sh = $s.Shapes.Item("Rounded Rectangle 35")
$sh.Left = 393.6961417322835
$sh.Top = 322.67984251968505

$sh = $s.Shapes.Item("TextBox 36")
$sh.Left = 444.6251968503937
$sh.Top = 407.4536220472441

$sh = $s.Shapes.Item("Oval 37")
$sh.Left = 414.8436220472441
$sh.Top = 340.2655118110236

$sh = $s.Shapes.Item("Oval 38")
$sh.Left = 473.09417322834645
$sh.Top = 328.4030708661417

$sh = $s.Shapes.Item("Oval 39")
$sh.Left = 489.15220472440944
$sh.Top = 371.1077165354331

$sh = $s.Shapes.Item("Straight Connector 40")
$sh.Left = 434.91614173228345
$sh.Top = 340.2655118110236

$sh = $s.Shapes.Item("Straight Arrow Connector 41")
$sh.Left = 483.13047244094486
$sh.Top = 352.12787401574803

$sh = $s.Shapes.Item("Straight Arrow Connector 42")
$sh.Left = 431.97661417322837
$sh.Top = 360.515905511811

$sh = $s.Shapes.Item("Curved Connector 43")
$sh.Left = 476.5536220472441
$sh.Top = 334.9799212598425

$sh = $s.Shapes.Item("Rectangle 44")
$sh.Left = 413.48818897637796
$sh.Top = 314.2240157480315

$sh = $s.Shapes.Item("TextBox 46")
$sh.Left = 354.28937007874015
$sh.Top = 293.7004724409449

$sh = $s.Shapes.Item("TextBox 47")
$sh.Left = 549.5535433070867
$sh.Top = 294.57314960629924

$sh = $s.Shapes.Item("Rectangle 48")
$sh.Left = 518.3041732283465
$sh.Top = 313.1899212598425

$sh = $s.Shapes.Item("Elbow Connector 54")
$sh.Left = 471.00976377952753
$sh.Top = 42.78314960629921

$sh = $s.Shapes.Item("Elbow Connector 56")
$sh.Left = 359.5126771653543
$sh.Top = 253.20299212598425

$sh = $s.Shapes.Item("Elbow Connector 61")
$sh.Left = 494.2203937007874
$sh.Top = 254.01629921259843

$sh = $s.Shapes.Item("Content Placeholder 2")
$sh.Left = 16.953543307086615
$sh.Top = 104.47850393700787
$sh.Width = 236.56165354330707
$sh.Height = 404.55417322834643
